$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 118
$ws.Range("I12").Value = 118
$ws.Range("K12").Value = 118
$ws.Range("M12").Value = 52

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 611.3
$ws.Range("J33").Value = 1640
$ws.Range("L33").Value = 1640
$ws.Range("N33").Value = -2098

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 719.4
$ws.Range("I38").Value = 719.4
$ws.Range("K38").Value = 2158.2
$ws.Range("M38").Value = -1786.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6966.16
$ws.Range("J64").Value = 8249.857
$ws.Range("L64").Value = 8249.857
$ws.Range("N64").Value = -8745.857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 6966.16
$ws.Range("J67").Value = 8249.857
$ws.Range("L67").Value = 8249.857
$ws.Range("N67").Value = -9965.857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2192.558
$ws.Range("J80").Value = 3413.68
$ws.Range("L80").Value = 10241.04
$ws.Range("N80").Value = -12237.04

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 2192.558
$ws.Range("J83").Value = 3413.68
$ws.Range("L83").Value = 30723.12
$ws.Range("N83").Value = -40707.12

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5550
$ws.Range("I86").Value = 4575
$ws.Range("K86").Value = 4575
$ws.Range("M86").Value = -3452

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 5550
$ws.Range("I89").Value = 4575
$ws.Range("K89").Value = 22875
$ws.Range("M89").Value = -17259

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 937.08
$ws.Range("I92").Value = 789.25
$ws.Range("J92").Value = 1199.8889
$ws.Range("K92").Value = 789.25
$ws.Range("L92").Value = 1199.8889
$ws.Range("M92").Value = 458.75
$ws.Range("N92").Value = -3695.8889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10755572
$ws.Range("I32").Value = 10872440
$ws.Range("K32").Value = 10872440
$ws.Range("M32").Value = -10872153

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2688.2222
$ws.Range("I45").Value = 2899.8
$ws.Range("K45").Value = 2899.8
$ws.Range("M45").Value = -2522.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1682.5385
$ws.Range("I97").Value = 1233.5758
$ws.Range("K97").Value = 1233.5758
$ws.Range("M97").Value = -737.5758000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 46999.855
$ws.Range("J138").Value = 46999.855
$ws.Range("L138").Value = 46999.855
$ws.Range("N138").Value = -57279.855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1064
$ws.Range("I64").Value = 925
$ws.Range("J64").Value = 1119.6
$ws.Range("K64").Value = 925
$ws.Range("L64").Value = 1119.6
$ws.Range("M64").Value = -700
$ws.Range("N64").Value = -1569.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 1064
$ws.Range("I67").Value = 925
$ws.Range("J67").Value = 1119.6
$ws.Range("K67").Value = 925
$ws.Range("L67").Value = 1119.6
$ws.Range("M67").Value = -145
$ws.Range("N67").Value = -2679.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4203.684
$ws.Range("I86").Value = 3428.4443
$ws.Range("J86").Value = 4901.4
$ws.Range("K86").Value = 3428.4443
$ws.Range("L86").Value = 4901.4
$ws.Range("M86").Value = -2305.4443
$ws.Range("N86").Value = -7147.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4203.684
$ws.Range("I89").Value = 3428.4443
$ws.Range("J89").Value = 4901.4
$ws.Range("K89").Value = 17142.2215
$ws.Range("L89").Value = 24507
$ws.Range("M89").Value = -11526.2215
$ws.Range("N89").Value = -35739

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1644.12
$ws.Range("I22").Value = 833
$ws.Range("J22").Value = 1754.7273
$ws.Range("K22").Value = 833
$ws.Range("L22").Value = 1754.7273
$ws.Range("M22").Value = -483
$ws.Range("N22").Value = -2454.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7914.923
$ws.Range("I62").Value = 8436.875
$ws.Range("J62").Value = 7079.8
$ws.Range("K62").Value = 8436.875
$ws.Range("L62").Value = 7079.8
$ws.Range("M62").Value = -7812.875
$ws.Range("N62").Value = -8327.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 7914.923
$ws.Range("I65").Value = 8436.875
$ws.Range("J65").Value = 7079.8
$ws.Range("K65").Value = 42184.375
$ws.Range("L65").Value = 35399
$ws.Range("M65").Value = -39064.375
$ws.Range("N65").Value = -41639

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 6662
$ws.Range("I105").Value = 6662
$ws.Range("K105").Value = 6662
$ws.Range("M105").Value = -4915

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1174.5
$ws.Range("I68").Value = 850
$ws.Range("J68").Value = 1282.6666
$ws.Range("K68").Value = 2550
$ws.Range("L68").Value = 3847.9998
$ws.Range("M68").Value = -1739
$ws.Range("N68").Value = -5469.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1174.5
$ws.Range("I71").Value = 850
$ws.Range("J71").Value = 1282.6666
$ws.Range("K71").Value = 7650
$ws.Range("L71").Value = 11543.9994
$ws.Range("M71").Value = -3594
$ws.Range("N71").Value = -19655.9994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 4498
$ws.Range("J106").Value = 4996
$ws.Range("L106").Value = 14988
$ws.Range("N106").Value = -16880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3216.2
$ws.Range("I137").Value = 1958.3334
$ws.Range("K137").Value = 5875.0002
$ws.Range("M137").Value = -775.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 64762.11
$ws.Range("I80").Value = 114687
$ws.Range("J80").Value = 2356
$ws.Range("K80").Value = 114687
$ws.Range("L80").Value = 2356
$ws.Range("M80").Value = -113689
$ws.Range("N80").Value = -4352

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 64762.11
$ws.Range("I83").Value = 114687
$ws.Range("J83").Value = 2356
$ws.Range("K83").Value = 573435
$ws.Range("L83").Value = 11780
$ws.Range("M83").Value = -568443
$ws.Range("N83").Value = -21764

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 8512.375
$ws.Range("I113").Value = 5025
$ws.Range("J113").Value = 11999.75
$ws.Range("K113").Value = 5025
$ws.Range("L113").Value = 11999.75
$ws.Range("M113").Value = -2855
$ws.Range("N113").Value = -16339.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 16504.857
$ws.Range("I40").Value = 18604.834
$ws.Range("K40").Value = 18604.834
$ws.Range("M40").Value = -18468.834

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3250.6
$ws.Range("I68").Value = 2588.25
$ws.Range("K68").Value = 2588.25
$ws.Range("M68").Value = -1839.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3250.6
$ws.Range("I71").Value = 2588.25
$ws.Range("K71").Value = 12941.25
$ws.Range("M71").Value = -9197.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5174.625
$ws.Range("I93").Value = 11750
$ws.Range("J93").Value = 2982.8333
$ws.Range("K93").Value = 11750
$ws.Range("L93").Value = 2982.8333
$ws.Range("M93").Value = -10502
$ws.Range("N93").Value = -5478.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 46731.668
$ws.Range("I100").Value = 77472.336
$ws.Range("J100").Value = 8305.833000000001
$ws.Range("K100").Value = 77472.336
$ws.Range("L100").Value = 8305.833000000001
$ws.Range("M100").Value = -76931.336
$ws.Range("N100").Value = -9387.833000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5142.1763

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 64000
$ws.Range("J46").Value = 64000
$ws.Range("L46").Value = 64000
$ws.Range("N46").Value = -64462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5803.4443
$ws.Range("J81").Value = 9333
$ws.Range("L81").Value = 18666
$ws.Range("N81").Value = -20788

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 5803.4443
$ws.Range("J84").Value = 9333
$ws.Range("L84").Value = 93330
$ws.Range("N84").Value = -103938

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2597.25
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 2796.3333
$ws.Range("K113").Value = 6000
$ws.Range("L113").Value = 8388.999899999999
$ws.Range("M113").Value = -3830
$ws.Range("N113").Value = -12728.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2383.6667
$ws.Range("I132").Value = 2241.0588
$ws.Range("J132").Value = 2989.75
$ws.Range("K132").Value = 6723.176399999999
$ws.Range("L132").Value = 8969.25
$ws.Range("M132").Value = -4193.176399999999
$ws.Range("N132").Value = -14029.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 64000
$ws.Range("J134").Value = 64000
$ws.Range("L134").Value = 192000
$ws.Range("N134").Value = -197070
